$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.393.57"
$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("D3").Value = "2.997.83"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.62%  "
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "2.984.29"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.37%  "
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "3.492.95"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.93%  "
$ws.Range("D18").Value = "2.996.41"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "59.341.16"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.722"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.83%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.83%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.12%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("D36").Value = "0.0₃0765"
$ws.Range("E36").Value = "  +9.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "406.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.37%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0354"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.776.28"
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("E51").Value = "  -0.19%  "
